$d = $word.ActiveDocument

# --- Replacement 1 (original line 142): intro sentence before the chatbot info ---
# The replacement text naturally ends with a space, so a plain Find/Replace
# keeps this run's xml:space="preserve" attribute (it is needed anyway).
$result0 = $d.Content.Find.Execute("Before you decide if you’d like to join, it’s important for you to know why we’re doing this research and what it involves. شما می‌توانید این ورقۀ معلومات اشتراک‌کننده را مطالعه نمایید. ", $true, $false, $false, $false, $false, $true, 1, $false, "پیش از آن‌که تصمیم بگیرید آیا می‌خواهید شامل این برنامه شوید یا خیر، مهم است بدانید که هدف از انجام این ارزیابی چیست و شامل چه مواردی می‌باشد.شما می‌توانید این ورقۀ معلومات اشتراک‌کننده را مطالعه نمایید. ", 2)
if (-not $result0) {
    foreach ($para in $d.Paragraphs) {
        if ($para.Range.Text -eq "Before you decide if you’d like to join, it’s important for you to know why we’re doing this research and what it involves. شما می‌توانید این ورقۀ معلومات اشتراک‌کننده را مطالعه نمایید. `r") {
            $rng0 = $para.Range.Duplicate()
            $rng0.Text = "پیش از آن‌که تصمیم بگیرید آیا می‌خواهید شامل این برنامه شوید یا خیر، مهم است بدانید که هدف از انجام این ارزیابی چیست و شامل چه مواردی می‌باشد.شما می‌توانید این ورقۀ معلومات اشتراک‌کننده را مطالعه نمایید. "
            $result0 = $true
            break
        }
    }
}
Write-Output "Replace 0: $result0"

# --- Replacement 2 (original line 640): "Data protection" heading ---
# The replacement text has no surrounding whitespace, so we replace the
# whole (single-run) paragraph Range directly (a duplicate of the
# untouched paragraph range), which this runtime serializes while still
# keeping xml:space="preserve" on the run, matching the source document.
$found1 = $false
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -eq "Data protection`r") {
        $rng1 = $para.Range.Duplicate()
        $rng1.Text = "حفاظت از معلومات"
        $found1 = $true
        break
    }
}
if (-not $found1) {
    $found1 = $d.Content.Find.Execute("Data protection", $true, $false, $false, $false, $false, $true, 1, $false, "حفاظت از معلومات", 2)
}
Write-Output "Replace 1: $found1"

# --- Replacement 3 (original line 653): data protection body paragraph ---
# This paragraph also contains an (unchanged) hyperlink run after this
# text, so we use Find/Replace which only touches the matched run. The
# replacement text ends with a space, so xml:space="preserve" is kept.
$result2 = $d.Content.Find.Execute("The University of Oxford, IDEMS, PLH, and World Vision make sure your information is used safely and correctly, just for research. The study follows data protection laws like GDPR (General Data Protection Regulation) in the UK and POPIA (Protection of Personal Information Act) in South Africa. You can learn more about your rights regarding your data by following this link: ", $true, $false, $false, $false, $false, $true, 1, $false, "دانشگاه آکسفورد، IDEMS، PLH و ورلد ویژن اطمینان می‌دهند که معلومات شما تنها برای اهداف پژوهشی و به‌گونه‌ی مصون و درست استفاده می‌گردد. این ارزیابی مطابق با قوانین حفاظت از داده‌ها، از جمله مقررهٔ عمومی حفاظت از داده‌ها (GDPR) در بریتانیا و قانون حفاظت از معلومات شخصی (POPIA) در افریقای جنوبی، اجرا می‌گردد. You can learn more about your rights regarding your data by following this link: ", 2)
Write-Output "Replace 2: $result2"

# --- Replacement 4 (original line 687): "Who has approved this study?" heading ---
$found3 = $false
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -eq "Who has approved this study?`r") {
        $rng3 = $para.Range.Duplicate()
        $rng3.Text = "این ارزیابی توسط کی ها تایید شده است؟"
        $found3 = $true
        break
    }
}
if (-not $found3) {
    $found3 = $d.Content.Find.Execute("Who has approved this study?", $true, $false, $false, $false, $false, $true, 1, $false, "این ارزیابی توسط کی ها تایید شده است؟", 2)
}
Write-Output "Replace 3: $found3"

# --- Replacement 5 (original line 700): approval committee body paragraph ---
$found4 = $false
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -eq "This study has received approval from the University of Oxford’s Social Sciences and Humanities Interdivisional Research Ethics Committee.`r") {
        $rng4 = $para.Range.Duplicate()
        $rng4.Text = "این ارزیابی از سوی کمیتهٔ اخلاق پژوهشی میان‌دانشکده‌ای علوم اجتماعی و علوم انسانی دانشگاه آکسفورد مورد تأیید قرار گرفته است."
        $found4 = $true
        break
    }
}
if (-not $found4) {
    $found4 = $d.Content.Find.Execute("This study has received approval from the University of Oxford’s Social Sciences and Humanities Interdivisional Research Ethics Committee.", $true, $false, $false, $false, $false, $true, 1, $false, "این ارزیابی از سوی کمیتهٔ اخلاق پژوهشی میان‌دانشکده‌ای علوم اجتماعی و علوم انسانی دانشگاه آکسفورد مورد تأیید قرار گرفته است.", 2)
}
Write-Output "Replace 4: $found4"
